$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source refreshed this monthly index table: each year's Oct/Nov/Dec rows now
# come first (ahead of that year's Jan-Sep rows), and 2022 + 2023(H1) data was
# appended, growing the sheet from 48 to 67 data rows (A2:D49 -> A2:D68).

# Copy the date-column cell format (border + bold + centered, style index 1 in the
# original file) so newly-created rows 50-68 in column A pick up the same look as
# the existing dated rows, instead of Excel fabricating a brand-new style entry.
$ws.Range("A2").Copy()
$ws.Range("A50:A68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
  @("2018-10", 100.4, 101.6, 100.5),
  @("2018-11", 100.6, 101.8, 100.8),
  @("2018-12", 101, 101.9, 100.8),
  @("2018-01", 99.76049999999999, 101.0406, 100.2468),
  @("2018-02", 99, 100.9, 100.1),
  @("2018-03", 98.7, 100.7, 99.59999999999999),
  @("2018-04", 98.5, 100.9, 99.09999999999999),
  @("2018-05", 98.90000000000001, 101, 99.5),
  @("2018-06", 99.2, 100.8, 99.8),
  @("2018-07", 100, 101.3, 100),
  @("2018-08", 100.7, 101.5, 100.4),
  @("2018-09", 100.7, 101.6, 100.4),
  @("2019-10", 100.3, 100.1, 100.2),
  @("2019-11", 99.3, 99.7, 99.90000000000001),
  @("2019-12", 98.8, 99.3, 99.8),
  @("2019-01", 100.9, 101.5, 100.9),
  @("2019-02", 101.3, 101.3, 100.5),
  @("2019-03", 101.3, 101.4, 100.6),
  @("2019-04", 101.2, 101.1, 101.3),
  @("2019-05", 101, 101, 101),
  @("2019-06", 101.1, 100.9, 100.6),
  @("2019-07", 100.3, 100.7, 100.3),
  @("2019-08", 99.8, 100.4, 100.5),
  @("2019-09", 100.1, 100.2, 100.7),
  @("2020-10", 96.90000000000001, 98.40000000000001, 98.7),
  @("2020-11", 97, 98.40000000000001, 98.59999999999999),
  @("2020-12", 97.2, 98.5, 98.8),
  @("2020-01", 99.3, 99.59999999999999, 99.3),
  @("2020-02", 99.5, 99.5, 99.40000000000001),
  @("2020-03", 98.90000000000001, 99, 99.59999999999999),
  @("2020-04", 99, 98.90000000000001, 99.59999999999999),
  @("2020-05", 99.09999999999999, 98.90000000000001, 99.8),
  @("2020-06", 99.09999999999999, 98.90000000000001, 99.7),
  @("2020-07", 99.2, 98.5, 99.59999999999999),
  @("2020-08", 98.59999999999999, 98.40000000000001, 99.09999999999999),
  @("2020-09", 98, 98.3, 98.59999999999999),
  @("2021-10", 102.3, 100.4, 101.2),
  @("2021-11", 102.8, 100.9, 100.8),
  @("2021-12", 102.6, 100.3, 100.9),
  @("2021-01", 98.7, 98.2, 100.1),
  @("2021-02", 98.59999999999999, 98.5, 99.90000000000001),
  @("2021-03", 100.1, 99.3, 100.2),
  @("2021-04", 100.3, 99.5, 100.3),
  @("2021-05", 100.6, 99.5, 100),
  @("2021-06", 100.3, 99.5, 99.59999999999999),
  @("2021-07", 100.9, 99.7, 99.90000000000001),
  @("2021-08", 101.6, 99.7, 100.5),
  @("2021-09", 101.7, 100, 100.7),
  @("2022-10", 104.6, 100.9, 102.7),
  @("2022-11", 104.3, 100.3, 103),
  @("2022-12", 104.5, 100.5, 103.1),
  @("2022-01", 103.1, 100.7, 101.6),
  @("2022-02", 103.4, 100.3, 102),
  @("2022-03", 102.8, 99.8, 101.4),
  @("2022-04", 102, 99.8, 101.3),
  @("2022-05", 102.9, 100.3, 101.6),
  @("2022-06", 103.6, 101.2, 102.5),
  @("2022-07", 103.5, 101.1, 102.7),
  @("2022-08", 103.4, 101.5, 102.4),
  @("2022-09", 104.1, 101.3, 102.4),
  @("2023-01", 103.6, 100.2, 102.4),
  @("2023-02", 103.3, 100.4, 101.8),
  @("2023-03", 103.6, 100.8, 102),
  @("2023-04", 103.4, 101.1, 101.7),
  @("2023-05", 102.4, 100.7, 101.4),
  @("2023-06", 102.1, 100.2, 101.3),
  @("2023-07", 102, 101.1, 100.9)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
